# The workbook's single data table (rows 2-269) contains one daily price
# record per row for "Zanahoria" at "Macroferia Regional de Talca".
# This edit inserts a new record (a duplicate of the data that was in row
# 141) right after row 141, shifting all subsequent rows down by one, and
# bumps the date in (the now edited) row 141 forward by one day
# (44586 -> 44587, i.e. 2022-02-18 -> 2022-02-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 142; rows 142-269 shift to 143-270.
$ws.Rows.Item(142).Insert()

# Duplicate row 141 (its original, still-unmodified content) into the
# newly inserted row 142.
$ws.Rows.Item(141).Copy()
$ws.Rows.Item(142).PasteSpecial()

# Now update the date of row 141 itself (44586 -> 44587).
$ws.Cells.Item(141, 4).Value2 = 44587
